# Workbook/worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before the current row 2. This pushes every
# existing song row (currently rows 2-257) down by 3 rows.
$ws.Rows.Item(2).Resize(3).Insert()

# The insert copies the formatting of row 1 (which has styled placeholder
# cells in columns A, D, E, F) into the new rows. Clear that so the new
# rows look like every other plain song row (only columns B and C used).
$ws.Rows.Item(2).Resize(3).Clear()

# Fill the three new rows at the top with the songs that now belong there:
#  - "One Direction - Story Of My Life" (previously further down the list)
#  - "Lewis Capaldi - Before You Go (Lyrics)" (brand new entry)
#  - "Harry Styles - Sign of the Times" (previously further down the list)
$ws.Range("B2").Value = "One Direction - Story Of My Life"
$ws.Range("C2").Value = "https://youtu.be/W-TE_Ys4iwM?si=RViOxRuaXxdz3pmm"

$ws.Range("B3").Value = "Lewis Capaldi - Before You Go (Lyrics)"
$ws.Range("C3").Value = "https://www.youtube.com/watch?v=c-KFmp9MMmQ"

$ws.Range("B4").Value = "Harry Styles - Sign of the Times"
$ws.Range("C4").Value = "https://www.youtube.com/watch?v=qN4ooNx77u0"

# Now remove the original occurrences of those two relocated songs, which
# (after the +3 row shift above) now live at row 166 ("Harry Styles - Sign
# of the Times", originally row 163) and row 68 ("One Direction - Story Of
# My Life", originally row 65). Delete the higher-numbered row first so the
# lower row number is unaffected by the shift caused by the deletion.
$ws.Rows.Item(166).Delete()
$ws.Rows.Item(68).Delete()
